# Fruta / hortaliza, semanal -- Femacal de La Calera (Papaya) weekly price refresh.
# Updates the date/volume/price fields for the existing daily rows and appends the
# new week (row 36), matching the upstream consolidated feed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column letter -> index, so edits below can be expressed as plain column letters.
$col = @{
    'A' = 1
    'B' = 2
    'C' = 3
    'D' = 4
    'E' = 5
    'F' = 6
    'G' = 7
    'H' = 8
    'I' = 9
    'J' = 10
    'K' = 11
    'L' = 12
    'M' = 13
    'N' = 14
    'O' = 15
    'P' = 16
    'Q' = 17
    'R' = 18
    'S' = 19
    'T' = 20
}

# Per-row edits: only the cells that actually change for that row.
$rowEdits = [ordered]@{
    2 = @{ 'D' = 44413; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'S' = 1500 }
    3 = @{ 'D' = 44413; 'M' = 58; 'N' = 13000; 'O' = 13000; 'P' = 13000; 'S' = 1300 }
    5 = @{ 'D' = 44315; 'M' = 60; 'N' = 24000; 'O' = 24000; 'P' = 24000; 'Q' = "`$/caja 15 kilos granel"; 'S' = 1600; 'T' = 15 }
    6 = @{ 'D' = 44454 }
    7 = @{ 'D' = 44398; 'M' = 60; 'N' = 17000; 'O' = 17000; 'P' = 17000; 'S' = 1700 }
    8 = @{ 'D' = 44398; 'L' = "Segunda"; 'M' = 50 }
    9 = @{ 'D' = 44396; 'L' = "Primera"; 'M' = 60; 'N' = 17000; 'O' = 17000; 'P' = 17000; 'S' = 1700 }
    10 = @{ 'D' = 44396; 'L' = "Segunda"; 'M' = 56; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = "`$/bandeja 10 kilos"; 'S' = 1500; 'T' = 10 }
    11 = @{ 'D' = 44329; 'M' = 50; 'N' = 16000; 'O' = 16000; 'P' = 16000; 'S' = 1600 }
    12 = @{ 'L' = "Primera"; 'M' = 65; 'N' = 18000; 'O' = 18000; 'P' = 18000; 'S' = 1800 }
    13 = @{ 'D' = 44431; 'L' = "Segunda"; 'M' = 60; 'N' = 16000; 'O' = 16000; 'P' = 16000; 'S' = 1600 }
    14 = @{ 'D' = 44435; 'M' = 115; 'N' = 18000; 'O' = 18000; 'P' = 18000; 'S' = 1800 }
    15 = @{ 'D' = 44435; 'M' = 60; 'N' = 16000; 'O' = 16000; 'P' = 16000; 'S' = 1600 }
    16 = @{ 'D' = 44445; 'M' = 68; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'S' = 1500 }
    17 = @{ 'D' = 44319; 'L' = "Primera"; 'M' = 60; 'N' = 24000; 'O' = 24000; 'P' = 24000; 'Q' = "`$/caja 15 kilos granel"; 'S' = 1600; 'T' = 15 }
    18 = @{ 'D' = 44385; 'M' = 60; 'N' = 17000; 'O' = 17000; 'P' = 17000; 'S' = 1700 }
    19 = @{ 'D' = 44385; 'L' = "Segunda"; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'S' = 1500 }
    20 = @{ 'D' = 44391; 'N' = 17000; 'O' = 17000; 'P' = 17000; 'S' = 1700 }
    21 = @{ 'D' = 44391; 'L' = "Segunda"; 'M' = 45; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = "`$/bandeja 10 kilos"; 'S' = 1500; 'T' = 10 }
    22 = @{ 'D' = 44417; 'M' = 56; 'N' = 16000; 'O' = 16000; 'P' = 16000; 'S' = 1600 }
    23 = @{ 'D' = 44417; 'M' = 60; 'N' = 14000; 'O' = 14000; 'P' = 14000; 'S' = 1400 }
    24 = @{ 'D' = 44321; 'M' = 42; 'N' = 24000; 'O' = 24000; 'P' = 24000; 'Q' = "`$/caja 15 kilos granel"; 'S' = 1600; 'T' = 15 }
    25 = @{ 'D' = 44441; 'L' = "Primera"; 'M' = 80 }
    26 = @{ 'D' = 44453; 'M' = 50; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'S' = 1500 }
    27 = @{ 'D' = 44354; 'L' = "Primera"; 'M' = 45; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'S' = 1500 }
    28 = @{ 'D' = 44420; 'M' = 54; 'N' = 18000; 'O' = 18000; 'P' = 18000; 'S' = 1800 }
    29 = @{ 'D' = 44420; 'L' = "Segunda"; 'M' = 50; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = "`$/bandeja 10 kilos"; 'S' = 1500; 'T' = 10 }
    30 = @{ 'D' = 44370 }
    31 = @{ 'D' = 44389; 'L' = "Primera"; 'M' = 60; 'N' = 17000; 'O' = 17000; 'P' = 17000; 'S' = 1700 }
    32 = @{ 'D' = 44389; 'L' = "Segunda"; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'S' = 1500 }
    33 = @{ 'D' = 44433; 'M' = 50 }
    34 = @{ 'D' = 44371; 'L' = "Primera"; 'M' = 67; 'N' = 18000; 'O' = 18000; 'P' = 18000; 'S' = 1800 }
    35 = @{ 'D' = 44410; 'M' = 75; 'N' = 15000; 'O' = 15000; 'P' = 15000; 'Q' = "`$/bandeja 10 kilos"; 'S' = 1500; 'T' = 10 }
}

foreach ($r in $rowEdits.Keys) {
    foreach ($c in $rowEdits[$r].Keys) {
        $ws.Cells.Item($r, $col[$c]).Value = $rowEdits[$r][$c]
    }
}

# New row 36: a full new weekly record appended after the last existing row.
$newRow = 36
$newRowValues = [ordered]@{
    'A' = 3
    'B' = "Femacal de La Calera"
    'C' = "Coquimbo"
    'D' = 44323
    'E' = 5
    'F' = "Fruta"
    'G' = 100108
    'H' = "Tropicales y subtropicales"
    'I' = 100108004
    'J' = "Papaya"
    'K' = "Cultivar IV Región"
    'L' = "Primera"
    'M' = 48
    'N' = 24000
    'O' = 24000
    'P' = 24000
    'Q' = "`$/caja 15 kilos granel"
    'R' = "Provincia del Elquí"
    'S' = 1600
    'T' = 15
}
foreach ($c in $newRowValues.Keys) {
    $ws.Cells.Item($newRow, $col[$c]).Value = $newRowValues[$c]
}

# Match the date formatting used by the other rows in column D.
$ws.Cells.Item($newRow, $col['D']).NumberFormat = "YYYY-MM-DD HH:MM:SS"

